$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 200
$ws.Range("I52").Value = 200
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 600
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -440
$ws.Range("N52").ClearContents()
$ws.Range("H109").Value = 32966
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 32966
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 32966
$ws.Range("N109").Value = -35740
$ws.Range("H112").Value = 1138001
$ws.Range("I112").Value = 1166.6666
$ws.Range("J112").Value = 1317501.2
$ws.Range("K112").Value = 3499.9998
$ws.Range("L112").Value = 3952503.6
$ws.Range("M112").Value = -2391.9998
$ws.Range("N112").Value = -3954719.6
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H132").Value = 33281.9
$ws.Range("I132").Value = 4915.5835
$ws.Range("J132").Value = 146747.17
$ws.Range("K132").Value = 14746.7505
$ws.Range("L132").Value = 440241.51
$ws.Range("M132").Value = -12216.7505
$ws.Range("N132").Value = -445301.51
$ws.Range("H137").Value = 5812.8623
$ws.Range("I137").Value = 1746.7333
$ws.Range("J137").Value = 10169.429
$ws.Range("K137").Value = 5240.199900000001
$ws.Range("L137").Value = 30508.287
$ws.Range("M137").Value = -2690.199900000001
$ws.Range("N137").Value = -35608.287
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11873.933
$ws.Range("I32").Value = 9456.672
$ws.Range("J32").Value = 23216.46
$ws.Range("K32").Value = 9456.672
$ws.Range("L32").Value = 23216.46
$ws.Range("M32").Value = -9169.672
$ws.Range("N32").Value = -23790.46
$ws.Range("H61").Value = 2678.5
$ws.Range("I61").Value = 1332.6875
$ws.Range("J61").Value = 3657.2727
$ws.Range("K61").Value = 1332.6875
$ws.Range("L61").Value = 3657.2727
$ws.Range("M61").Value = -1120.6875
$ws.Range("N61").Value = -4081.2727
$ws.Range("H111").Value = 48140
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 48140
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 48140
$ws.Range("N111").Value = -56320
$ws.Range("H114").Value = 43087.5
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 43087.5
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 43087.5
$ws.Range("N114").Value = -51765.5
$ws.Range("H117").Value = 38551.2
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 38551.2
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 38551.2
$ws.Range("N117").Value = -47729.2
$ws.Range("H118").Value = 49376
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 49376
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 49376
$ws.Range("N118").Value = -52690
$ws.Range("H136").Value = 2678.5
$ws.Range("I136").Value = 1332.6875
$ws.Range("J136").Value = 3657.2727
$ws.Range("K136").Value = 3998.0625
$ws.Range("L136").Value = 10971.8181
$ws.Range("M136").Value = -1448.0625
$ws.Range("N136").Value = -16071.8181
$ws.Range("H138").Value = 28311
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 28311
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 28311
$ws.Range("N138").Value = -38591
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H116").Value = 43415.5
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 43415.5
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 43415.5
$ws.Range("N116").Value = -52593.5
$ws.Range("H117").Value = 44517
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 44517
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 44517
$ws.Range("N117").Value = -53695
$ws.Range("H124").Value = 49322
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 49322
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 49322
$ws.Range("N124").Value = -59142
$ws.Range("H134").Value = 233962.62
$ws.Range("I134").Value = 3302.1667
$ws.Range("J134").Value = 259125.58
$ws.Range("K134").Value = 9906.500100000001
$ws.Range("L134").Value = 777376.74
$ws.Range("M134").Value = -7371.500100000001
$ws.Range("N134").Value = -782446.74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1859.963
$ws.Range("I58").Value = 1780.1111
$ws.Range("J58").Value = 1899.8889
$ws.Range("K58").Value = 1780.1111
$ws.Range("L58").Value = 1899.8889
$ws.Range("M58").Value = -1577.1111
$ws.Range("N58").Value = -2305.8889
$ws.Range("H100").Value = 46964.25
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 46964.25
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 46964.25
$ws.Range("N100").Value = -49128.25
$ws.Range("H116").Value = 42454
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 42454
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 42454
$ws.Range("N116").Value = -51632
$ws.Range("H134").Value = 484787.9
$ws.Range("I134").Value = 1473.2941
$ws.Range("J134").Value = 1169483.6
$ws.Range("K134").Value = 4419.8823
$ws.Range("L134").Value = 3508450.8
$ws.Range("M134").Value = -1884.8823
$ws.Range("N134").Value = -3513520.8
$ws.Range("H136").Value = 1859.963
$ws.Range("I136").Value = 1780.1111
$ws.Range("J136").Value = 1899.8889
$ws.Range("K136").Value = 5340.3333
$ws.Range("L136").Value = 5699.6667
$ws.Range("M136").Value = -2790.3333
$ws.Range("N136").Value = -10799.6667
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3489.4856
$ws.Range("I113").Value = 7316.2666
$ws.Range("J113").Value = 619.4
$ws.Range("K113").Value = 21948.7998
$ws.Range("L113").Value = 1858.2
$ws.Range("M113").Value = -19778.7998
$ws.Range("N113").Value = -6198.2
$ws.Range("H136").Value = 50002356
$ws.Range("I136").Value = 71430240
$ws.Range("J136").Value = 3955.6667
$ws.Range("K136").Value = 214290720
$ws.Range("L136").Value = 11867.0001
$ws.Range("M136").Value = -214285620
$ws.Range("N136").Value = -22067.0001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 45640
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 45640
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 45640
$ws.Range("N130").Value = -55680
$ws.Range("H132").Value = 3033.457
$ws.Range("I132").Value = 1967.6316
$ws.Range("J132").Value = 4299.125
$ws.Range("K132").Value = 5902.8948
$ws.Range("L132").Value = 12897.375
$ws.Range("M132").Value = -3372.8948
$ws.Range("N132").Value = -17957.375
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 48571
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 48571
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 48571
$ws.Range("N108").Value = -56251
$ws.Range("H110").Value = 30000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 30000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 30000
$ws.Range("N110").Value = -38180
$ws.Range("H132").Value = 2529.4314
$ws.Range("I132").Value = 1729.9412
$ws.Range("J132").Value = 4128.4116
$ws.Range("K132").Value = 5189.8236
$ws.Range("L132").Value = 12385.2348
$ws.Range("M132").Value = -2659.8236
$ws.Range("N132").Value = -17445.2348
$ws.Range("H136").Value = 1332.0344
$ws.Range("I136").Value = 1023.32434
$ws.Range("J136").Value = 1875.9524
$ws.Range("K136").Value = 3069.97302
$ws.Range("L136").Value = 5627.857199999999
$ws.Range("M136").Value = -519.9730199999999
$ws.Range("N136").Value = -10727.8572
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 45984
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 45984
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 45984
$ws.Range("N117").Value = -55162
$ws.Range("H132").Value = 2853.2
$ws.Range("I132").Value = 1643.2858
$ws.Range("J132").Value = 3911.875
$ws.Range("K132").Value = 4929.857400000001
$ws.Range("L132").Value = 11735.625
$ws.Range("M132").Value = -2399.857400000001
$ws.Range("N132").Value = -16795.625
